{"js": "// Replace the date string and each two-digit multiplication equation with\n// the updated values, per the commit's diff. Every old string is unique in\n// the document, so a literal (non-wildcard) search/replace per pair is safe.\n\nconst replacements = [\n  [\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"],\n  [\"88\u00d779=6952\", \"37\u00d795=3515\"],\n  [\"91\u00d757=5187\", \"73\u00d768=4964\"],\n  [\"56\u00d727=1512\", \"68\u00d743=2924\"],\n  [\"33\u00d779=2607\", \"17\u00d775=1275\"],\n  [\"34\u00d732=1088\", \"77\u00d774=5698\"],\n  [\"52\u00d754=2808\", \"15\u00d789=1335\"],\n  [\"74\u00d744=3256\", \"48\u00d746=2208\"],\n  [\"93\u00d716=1488\", \"27\u00d787=2349\"],\n  [\"63\u00d721=1323\", \"75\u00d711=825\"],\n  [\"26\u00d769=1794\", \"66\u00d760=3960\"],\n  [\"18\u00d756=1008\", \"91\u00d798=8918\"],\n  [\"23\u00d751=1173\", \"82\u00d740=3280\"],\n  [\"93\u00d723=2139\", \"56\u00d780=4480\"],\n  [\"69\u00d792=6348\", \"81\u00d784=6804\"],\n  [\"19\u00d738=722\", \"99\u00d735=3465\"],\n  [\"48\u00d724=1152\", \"79\u00d742=3318\"],\n  [\"97\u00d771=6887\", \"72\u00d780=5760\"],\n  [\"33\u00d796=3168\", \"28\u00d745=1260\"],\n  [\"37\u00d744=1628\", \"53\u00d760=3180\"],\n  [\"32\u00d718=576\", \"62\u00d737=2294\"],\n  [\"22\u00d724=528\", \"24\u00d780=1920\"],\n  [\"51\u00d794=4794\", \"57\u00d797=5529\"],\n  [\"86\u00d791=7826\", \"74\u00d752=3848\"],\n  [\"45\u00d728=1260\", \"17\u00d730=510\"],\n  [\"12\u00d768=816\", \"82\u00d743=3526\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date string and each two-digit multiplication equation with\n# the updated values, per the commit's diff. Every old string is unique in\n# the document, so a literal Find/Replace (no wildcards) per pair is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"),\n    @(\"88\u00d779=6952\", \"37\u00d795=3515\"),\n    @(\"91\u00d757=5187\", \"73\u00d768=4964\"),\n    @(\"56\u00d727=1512\", \"68\u00d743=2924\"),\n    @(\"33\u00d779=2607\", \"17\u00d775=1275\"),\n    @(\"34\u00d732=1088\", \"77\u00d774=5698\"),\n    @(\"52\u00d754=2808\", \"15\u00d789=1335\"),\n    @(\"74\u00d744=3256\", \"48\u00d746=2208\"),\n    @(\"93\u00d716=1488\", \"27\u00d787=2349\"),\n    @(\"63\u00d721=1323\", \"75\u00d711=825\"),\n    @(\"26\u00d769=1794\", \"66\u00d760=3960\"),\n    @(\"18\u00d756=1008\", \"91\u00d798=8918\"),\n    @(\"23\u00d751=1173\", \"82\u00d740=3280\"),\n    @(\"93\u00d723=2139\", \"56\u00d780=4480\"),\n    @(\"69\u00d792=6348\", \"81\u00d784=6804\"),\n    @(\"19\u00d738=722\", \"99\u00d735=3465\"),\n    @(\"48\u00d724=1152\", \"79\u00d742=3318\"),\n    @(\"97\u00d771=6887\", \"72\u00d780=5760\"),\n    @(\"33\u00d796=3168\", \"28\u00d745=1260\"),\n    @(\"37\u00d744=1628\", \"53\u00d760=3180\"),\n    @(\"32\u00d718=576\", \"62\u00d737=2294\"),\n    @(\"22\u00d724=528\", \"24\u00d780=1920\"),\n    @(\"51\u00d794=4794\", \"57\u00d797=5529\"),\n    @(\"86\u00d791=7826\", \"74\u00d752=3848\"),\n    @(\"45\u00d728=1260\", \"17\u00d730=510\"),\n    @(\"12\u00d768=816\", \"82\u00d743=3526\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
